# Changes in Payment Methods
# - Adds three new "Payment Methods" test-case rows to Sheet1
# - Changes B31 Execute flag from Yes to No
# - Clears the now-unused helper table (F17:G28) on the "lists" sheet

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("lists")

# --- Sheet1: flip Execute flag for the "Verify Forgot Pin With Navigation" row ---
$ws1.Range("B31").Value = "No"

# --- Sheet1: add the three new Payment Methods rows, copying formatting from a
#     similar existing row (row 25) so fonts / alignment / number formats match. ---
$ws1.Range("A25:I25").Copy()
$ws1.Range("A32:I34").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Rows.Item(32).RowHeight = 172.8
$ws1.Rows.Item(33).RowHeight = 144
$ws1.Rows.Item(34).RowHeight = 79.2

# Row 32 - Verify Add Debit Card
$ws1.Range("A32").Value = "Verify Add Debit Card"
$ws1.Range("B32").Value = "No"
$ws1.Range("C32").Value = "testdata.xls,PaymentMethods"
$ws1.Range("D32").Value = "RunOneIteration"
$ws1.Range("E32").Value = "'1"
$ws1.Range("F32").Value = "'1"
$ws1.Range("G32").Value = "PaymentMethods"
$ws1.Range("H32").Value = "coyni_mobile.tests.LoginTest,`ntestLogin,`n-pemail,`n-ppassword,`n-ppin,`n-puserName"
$ws1.Range("I32").Value = "coyni_mobile.tests.CustomerProfileTest,`ntestAddDebitCard,`n-pnameOnCard,`n-pcardNumber,`n-pcardExp,`n-pcvvOrCVC,`n-paddressLine1,`n-paddreddLine2,`n-pcity,`n-pstate,`n-pzipCode,`n-pamount"

# Row 33 - Verify Edit Debit Card In Payment Methods
$ws1.Range("A33").Value = "Verify Edit Debit Card In Payment Methods"
$ws1.Range("B33").Value = "No"
$ws1.Range("C33").Value = "testdata.xls,PaymentMethods"
$ws1.Range("D33").Value = "RunOneIteration"
$ws1.Range("E33").Value = "'1"
$ws1.Range("F33").Value = "'1"
$ws1.Range("G33").Value = "PaymentMethods"
$ws1.Range("H33").Value = "coyni_mobile.tests.LoginTest,`ntestLogin,`n-pemail,`n-ppassword,`n-ppin,`n-puserName"
$ws1.Range("I33").Value = "coyni_mobile.tests.CustomerProfileTest,`ntestEditDebitCard,`n-pcardNumber,`n-pcardExp,`n-pcvvOrCVC,`n-paddressLine1,`n-paddreddLine2,`n-pcity,`n-pstate,`n-pzipCode"

# Row 34 - Verify Delete Debit Card In Payment Methods
$ws1.Range("A34").Value = "Verify Delete Debit Card In Payment Methods"
$ws1.Range("B34").Value = "Yes"
$ws1.Range("C34").Value = "testdata.xls,PaymentMethods"
$ws1.Range("D34").Value = "RunOneIteration"
$ws1.Range("E34").Value = "'1"
$ws1.Range("F34").Value = "'1"
$ws1.Range("G34").Value = "PaymentMethods"
$ws1.Range("H34").Value = "coyni_mobile.tests.LoginTest,`ntestLogin,`n-pemail,`n-ppassword,`n-ppin,`n-puserName"
$ws1.Range("I34").Value = "coyni_mobile.tests.CustomerProfileTest,`ntestDeleteDebitCard,`n-pcardNumber"

# --- lists sheet: remove the stale helper table that lived in F17:G28 ---
$ws2.Range("F17:G28").Clear()

# --- Restore view/selection state on both sheets ---
$ws2.Activate()
$ws2.Range("G19").Select()

$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 2
$ws1.Range("H43").Select()
